$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 11049.333
$ws.Range("I62").Value = 14107.883
$ws.Range("K62").Value = 14107.883
$ws.Range("M62").Value = -13483.883
# Row 65
$ws.Range("H65").Value = 11049.333
$ws.Range("I65").Value = 14107.883
$ws.Range("K65").Value = 70539.41499999999
$ws.Range("M65").Value = -67419.41499999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 35
$ws.Range("H35").Value = 3603.4
$ws.Range("I35").Value = 2504.25
$ws.Range("J35").Value = 8000
$ws.Range("K35").Value = 2504.25
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = -2098.25
$ws.Range("N35").Value = -8812
# Row 45
$ws.Range("H45").Value = 1378
$ws.Range("I45").Value = 1106
$ws.Range("J45").Value = 1595.6
$ws.Range("K45").Value = 1106
$ws.Range("L45").Value = 1595.6
$ws.Range("M45").Value = -729
$ws.Range("N45").Value = -2349.6
# Row 88
$ws.Range("H88").Value = 19430.684
$ws.Range("I88").Value = 24757.45
$ws.Range("J88").Value = 2266.6667
$ws.Range("K88").Value = 24757.45
$ws.Range("L88").Value = 2266.6667
$ws.Range("M88").Value = -24351.45
$ws.Range("N88").Value = -3078.6667
# Row 91
$ws.Range("H91").Value = 19430.684
$ws.Range("I91").Value = 24757.45
$ws.Range("J91").Value = 2266.6667
$ws.Range("K91").Value = 24757.45
$ws.Range("L91").Value = 2266.6667
$ws.Range("M91").Value = -23353.45
$ws.Range("N91").Value = -5074.6667
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
# Row 134
$ws.Range("H134").Value = 29993.334
$ws.Range("J134").Value = 29993.334
$ws.Range("L134").Value = 29993.334
$ws.Range("N134").Value = -40133.334

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 37
$ws.Range("H37").Value = 980.2
$ws.Range("I37").Value = 475.25
$ws.Range("K37").Value = 475.25
$ws.Range("M37").Value = -338.25
# Row 86
$ws.Range("H86").Value = 4848.476
$ws.Range("I86").Value = 5471.2354
$ws.Range("J86").Value = 2201.75
$ws.Range("K86").Value = 5471.2354
$ws.Range("L86").Value = 2201.75
$ws.Range("M86").Value = -4348.2354
$ws.Range("N86").Value = -4447.75
# Row 89
$ws.Range("H89").Value = 4848.476
$ws.Range("I89").Value = 5471.2354
$ws.Range("J89").Value = 2201.75
$ws.Range("K89").Value = 27356.177
$ws.Range("L89").Value = 11008.75
$ws.Range("M89").Value = -21740.177
$ws.Range("N89").Value = -22240.75
# Row 138
$ws.Range("H138").Value = 26500
$ws.Range("J138").Value = 26500
$ws.Range("L138").Value = 26500
$ws.Range("N138").Value = -36780

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 959187.0600000001
$ws.Range("I31").Value = 793.7083
$ws.Range("J31").Value = 1937971.8
$ws.Range("K31").Value = 793.7083
$ws.Range("L31").Value = 1937971.8
$ws.Range("M31").Value = -498.7083
$ws.Range("N31").Value = -1938561.8
# Row 34
$ws.Range("H34").Value = 959187.0600000001
$ws.Range("I34").Value = 793.7083
$ws.Range("J34").Value = 1937971.8
$ws.Range("K34").Value = 793.7083
$ws.Range("L34").Value = 1937971.8
$ws.Range("M34").Value = -591.7083
$ws.Range("N34").Value = -1938375.8
# Row 106
$ws.Range("H106").Value = 49000
$ws.Range("J106").Value = 49000
$ws.Range("L106").Value = 49000
$ws.Range("N106").Value = -51524

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 3557.1428
$ws.Range("I25").Value = 490
$ws.Range("J25").Value = 4784
$ws.Range("K25").Value = 1470
$ws.Range("L25").Value = 14352
$ws.Range("M25").Value = -1301
$ws.Range("N25").Value = -14690
# Row 30
$ws.Range("H30").Value = 3557.1428
$ws.Range("I30").Value = 490
$ws.Range("J30").Value = 4784
$ws.Range("K30").Value = 1470
$ws.Range("L30").Value = 14352
$ws.Range("M30").Value = -1368
$ws.Range("N30").Value = -14556
# Row 34
$ws.Range("H34").Value = 935
$ws.Range("J34").Value = 935
$ws.Range("L34").Value = 2805
$ws.Range("N34").Value = -2973
# Row 39
$ws.Range("H39").Value = 4380
$ws.Range("J39").Value = 4380
$ws.Range("L39").Value = 13140
$ws.Range("N39").Value = -13728
# Row 40
$ws.Range("H40").Value = 150.39285
$ws.Range("I40").Value = 106.708336
$ws.Range("J40").Value = 412.5
$ws.Range("K40").Value = 426.833344
$ws.Range("L40").Value = 1650
$ws.Range("M40").Value = -357.833344
$ws.Range("N40").Value = -1788
# Row 55
$ws.Range("H55").Value = 3299.1667
$ws.Range("J55").Value = 3879
$ws.Range("L55").Value = 11637
$ws.Range("N55").Value = -11991
# Row 68
$ws.Range("H68").Value = 1195.279
$ws.Range("I68").Value = 629.0476
$ws.Range("J68").Value = 1735.7727
$ws.Range("K68").Value = 1887.1428
$ws.Range("L68").Value = 5207.3181
$ws.Range("M68").Value = -1076.1428
$ws.Range("N68").Value = -6829.3181
# Row 71
$ws.Range("H71").Value = 1195.279
$ws.Range("I71").Value = 629.0476
$ws.Range("J71").Value = 1735.7727
$ws.Range("K71").Value = 5661.4284
$ws.Range("L71").Value = 15621.9543
$ws.Range("M71").Value = -1605.4284
$ws.Range("N71").Value = -23733.9543
# Row 75
$ws.Range("H75").Value = 1301.8334
$ws.Range("I75").Value = 105.5
$ws.Range("J75").Value = 1900
$ws.Range("K75").Value = 316.5
$ws.Range("L75").Value = 5700
$ws.Range("M75").Value = 681.5
$ws.Range("N75").Value = -7696
# Row 78
$ws.Range("H78").Value = 1301.8334
$ws.Range("I78").Value = 105.5
$ws.Range("J78").Value = 1900
$ws.Range("K78").Value = 949.5
$ws.Range("L78").Value = 17100
$ws.Range("M78").Value = 4042.5
$ws.Range("N78").Value = -27084
# Row 102
$ws.Range("H102").Value = 8322.223
$ws.Range("I102").Value = 2450
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 7350
$ws.Range("L102").Value = 30000
$ws.Range("M102").Value = -4916
$ws.Range("N102").Value = -34868

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8371.071
$ws.Range("I70").Value = 5261
$ws.Range("J70").Value = 9047.174000000001
$ws.Range("K70").Value = 5261
$ws.Range("L70").Value = 9047.174000000001
$ws.Range("M70").Value = -4991
$ws.Range("N70").Value = -9587.174000000001
# Row 73
$ws.Range("H73").Value = 8371.071
$ws.Range("I73").Value = 5261
$ws.Range("J73").Value = 9047.174000000001
$ws.Range("K73").Value = 5261
$ws.Range("L73").Value = 9047.174000000001
$ws.Range("M73").Value = -4325
$ws.Range("N73").Value = -10919.174
# Row 102
$ws.Range("H102").Value = 4098.5713
$ws.Range("I102").Value = 2473.5386
$ws.Range("K102").Value = 2473.5386
$ws.Range("M102").Value = -851.5385999999999
# Row 105
$ws.Range("H105").Value = 34300
$ws.Range("J105").Value = 34300
$ws.Range("L105").Value = 34300
$ws.Range("N105").Value = -41288
# Row 122
$ws.Range("H122").Value = 1235.6923
$ws.Range("I122").Value = 1032.8334
$ws.Range("J122").Value = 1409.5714
$ws.Range("K122").Value = 3098.5002
$ws.Range("L122").Value = 4228.7142
$ws.Range("M122").Value = -648.5001999999999
$ws.Range("N122").Value = -9128.7142
# Row 126
$ws.Range("H126").Value = 2616.7144
$ws.Range("I126").Value = 2263.4
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 6790.200000000001
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -4320.200000000001
$ws.Range("N126").Value = -15440

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2367.7856
$ws.Range("I7").Value = 2545.2727
$ws.Range("J7").Value = 2252.9412
$ws.Range("K7").Value = 2545.2727
$ws.Range("L7").Value = 2252.9412
$ws.Range("M7").Value = -2433.2727
$ws.Range("N7").Value = -2476.9412
# Row 16
$ws.Range("H16").Value = 2266.1538
$ws.Range("I16").Value = 2469.0908
$ws.Range("J16").Value = 1150
$ws.Range("K16").Value = 2469.0908
$ws.Range("L16").Value = 1150
$ws.Range("M16").Value = -2299.0908
$ws.Range("N16").Value = -1490
# Row 40
$ws.Range("H40").Value = 2154.9092
$ws.Range("I40").Value = 2084
$ws.Range("J40").Value = 2240
$ws.Range("K40").Value = 2084
$ws.Range("L40").Value = 2240
$ws.Range("M40").Value = -1948
$ws.Range("N40").Value = -2512
# Row 46
$ws.Range("H46").Value = 1290.1464
$ws.Range("I46").Value = 1120.4073
$ws.Range("J46").Value = 1617.5
$ws.Range("K46").Value = 1120.4073
$ws.Range("L46").Value = 1617.5
$ws.Range("M46").Value = -932.4073000000001
$ws.Range("N46").Value = -1993.5
# Row 126
$ws.Range("H126").Value = 2367.7856
$ws.Range("I126").Value = 2545.2727
$ws.Range("J126").Value = 2252.9412
$ws.Range("K126").Value = 7635.8181
$ws.Range("L126").Value = 6758.823600000001
$ws.Range("M126").Value = -5165.8181
$ws.Range("N126").Value = -11698.8236

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("H104").Value = 32720
$ws.Range("J104").Value = 32720
$ws.Range("L104").Value = 32720
$ws.Range("N104").Value = -39708
# Row 126
$ws.Range("H126").Value = 2336.5
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 2418.3333
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 7254.3333
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -12194.9999
